# Scheduled-runner update: refresh computed profit figures across all
# recipe sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect the
# latest market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 625384.75
$ws.Range("J6").Value = 1504.5
$ws.Range("L6").Value = 4513.5
$ws.Range("N6").Value = -4737.5

$ws.Range("H28").Value = 2422.923
$ws.Range("I28").Value = 1462.375
$ws.Range("J28").Value = 3959.8
$ws.Range("K28").Value = 1462.375
$ws.Range("L28").Value = 3959.8
$ws.Range("M28").Value = -977.375
$ws.Range("N28").Value = -4929.8

$ws.Range("H33").Value = 126.916664
$ws.Range("I33").Value = 126.916664
$ws.Range("K33").Value = 126.916664
$ws.Range("M33").Value = 102.083336

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 4200
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 2000
$ws.Range("M64").Value = -1752

$ws.Range("H67").Value = 4200
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 2000
$ws.Range("M67").Value = -1142

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H113").Value = 5999.6665
$ws.Range("I113").Value = 3999.5
$ws.Range("K113").Value = 3999.5
$ws.Range("M113").Value = -745.5

$ws.Range("H138").Value = 9166.333
$ws.Range("J138").Value = 9750
$ws.Range("L138").Value = 29250
$ws.Range("N138").Value = -39530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 371
$ws.Range("I22").Value = 371
$ws.Range("K22").Value = 371
$ws.Range("M22").Value = -198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 178.75
$ws.Range("I2").Value = 178.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 178.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -65.75
$ws.Range("N2").ClearContents()

$ws.Range("H28").Value = 32999.5
$ws.Range("J28").Value = 32999.5
$ws.Range("L28").Value = 32999.5
$ws.Range("N28").Value = -33489.5

$ws.Range("H31").Value = 5509.2666
$ws.Range("I31").Value = 1701.2222
$ws.Range("K31").Value = 1701.2222
$ws.Range("M31").Value = -1406.2222

$ws.Range("H34").Value = 5509.2666
$ws.Range("I34").Value = 1701.2222
$ws.Range("K34").Value = 1701.2222
$ws.Range("M34").Value = -1499.2222

$ws.Range("H62").Value = 4995
$ws.Range("I62").Value = 4989
$ws.Range("K62").Value = 4989
$ws.Range("M62").Value = -4365

$ws.Range("H65").Value = 4995
$ws.Range("I65").Value = 4989
$ws.Range("K65").Value = 24945
$ws.Range("M65").Value = -21825

$ws.Range("H68").Value = 49800
$ws.Range("J68").Value = 49800
$ws.Range("L68").Value = 49800
$ws.Range("N68").Value = -51298

$ws.Range("H71").Value = 49800
$ws.Range("J71").Value = 49800
$ws.Range("L71").Value = 149400
$ws.Range("N71").Value = -156888

$ws.Range("H74").Value = 27165.834
$ws.Range("J74").Value = 29999
$ws.Range("L74").Value = 29999
$ws.Range("N74").Value = -31747

$ws.Range("H77").Value = 27165.834
$ws.Range("J77").Value = 29999
$ws.Range("L77").Value = 89997
$ws.Range("N77").Value = -98733

$ws.Range("H106").Value = 49998.25
$ws.Range("J106").Value = 49998.25
$ws.Range("L106").Value = 49998.25
$ws.Range("N106").Value = -52522.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 834.3333
$ws.Range("I8").Value = 834.3333
$ws.Range("K8").Value = 2502.9999
$ws.Range("M8").Value = -2363.9999

$ws.Range("H34").Value = 314
$ws.Range("I34").Value = 228
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 684
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -600
$ws.Range("N34").Value = -1368

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354

$ws.Range("H134").Value = 4486.857
$ws.Range("I134").Value = 4486.857
$ws.Range("K134").Value = 13460.571
$ws.Range("M134").Value = -8390.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62
$ws.Range("I2").Value = 61.666668
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 61.666668
$ws.Range("L2").Value = 63
$ws.Range("M2").Value = 51.333332
$ws.Range("N2").Value = -289

$ws.Range("H18").Value = 1000000
$ws.Range("I18").Value = 1000000
$ws.Range("K18").Value = 1000000
$ws.Range("M18").Value = -999707

$ws.Range("H46").Value = 490.625
$ws.Range("I46").Value = 490.625
$ws.Range("K46").Value = 490.625
$ws.Range("M46").Value = -334.625

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 699.75
$ws.Range("I22").Value = 399.33334
$ws.Range("K22").Value = 399.33334
$ws.Range("M22").Value = -104.33334

$ws.Range("H27").Value = 699.75
$ws.Range("I27").Value = 399.33334
$ws.Range("K27").Value = 399.33334
$ws.Range("M27").Value = -292.33334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 314.4
$ws.Range("I107").Value = 350.625
$ws.Range("J107").Value = 169.5
$ws.Range("K107").Value = 1051.875
$ws.Range("L107").Value = 508.5
$ws.Range("M107").Value = 868.125
$ws.Range("N107").Value = -4348.5

$ws.Range("H135").Value = 74998
$ws.Range("J135").Value = 74998
$ws.Range("L135").Value = 74998
$ws.Range("N135").Value = -85138

$ws.Range("H141").Value = 99927
$ws.Range("J141").Value = 99927
$ws.Range("L141").Value = 99927
$ws.Range("N141").Value = -110287

